# Rename the four header/footer logo pictures (Word "Name" property —
# stored as the wp:docPr/@name, mirrored by pic:cNvPr/@name in the part's
# markup). The Pearson logo (footer, both default + first-page) goes from
# image2.png -> image1.png; the BTec logo (header, both default + first-
# page) goes from image1.jpg -> image2.jpg.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

$footers = $sec.Footers
$headers = $sec.Headers

# Pearson logo, default footer
$pic = $footers.Item(1).Range.InlineShapes.Item(1)
$pic.Name = "image1.png"

# Pearson logo, first-page footer
$pic = $footers.Item(2).Range.InlineShapes.Item(1)
$pic.Name = "image1.png"

# BTec logo, default header
$pic = $headers.Item(1).Range.InlineShapes.Item(1)
$pic.Name = "image2.jpg"

# BTec logo, first-page header
$pic = $headers.Item(2).Range.InlineShapes.Item(1)
$pic.Name = "image2.jpg"
